# This workbook contains a weekly data series (Macroferia Regional de Talca -
# Brocoli). A new week's record is inserted at row 86 (the most recent
# record), which pushes every subsequent row's "variable" data (Fecha,
# Calidad, Volumen, Precio minimo/maximo/promedio ponderado, Origen,
# Precio $/Kg) down by one row, and the former last row (203) becomes a new
# row 204 at the bottom of the sheet.
#
# The columns that are constant across the whole data block (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) stay the same and simply
# need to be filled in for the new row 204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the variable columns for rows 86..203 down to rows 87..204 ---
# Column D (Fecha)
$srcD = $ws.Range("D86:D203")
$valsD = $srcD.Value()
$ws.Range("D87:D204").Value = $valsD

# Columns I..P (Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen, Precio $/Kg)
$srcIP = $ws.Range("I86:P203")
$valsIP = $srcIP.Value()
$ws.Range("I87:P204").Value = $valsIP

# --- Fill in the constant columns for the new row 204 (same values as the
# rest of the data block) ---
# Make sure the date cell uses the same number format as the rest of the
# "Fecha" column instead of Excel's auto-detected date format.
$ws.Range("D204").NumberFormat = $ws.Range("D203").NumberFormat()
$ws.Range("A204").Value = 5
$ws.Range("B204").Value = "Macroferia Regional de Talca"
$ws.Range("C204").Value = "Maule"
$ws.Range("E204").Value = 7
$ws.Range("F204").Value = 100112023
$ws.Range("G204").Value = "Brócoli"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"

# --- Write the brand-new record into row 86 ---
$ws.Range("D86").Value = 44482
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 800
$ws.Range("L86").Value = 800
$ws.Range("M86").Value = 800
$ws.Range("N86").Value = "$/unidad"
$ws.Range("O86").Value = "Región Metropolitana"
$ws.Range("P86").Value = 800
